$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: A2 -> "Gezinssamenstelling", add H2 = "x"
$ws.Range("A2").Value = "Gezinssamenstelling"
$ws.Range("H2").Value = "x"

# Update row 6: A6 -> "Eigen risico", add H6 = "x"
$ws.Range("A6").Value = "Eigen risico"
$ws.Range("H6").Value = "x"

# Update row 9: A9 -> "n/a", add H9 = "x"
$ws.Range("A9").Value = "n/a"
$ws.Range("H9").Value = "x"

# Remove rows 14-21 (the "asd" filler rows) entirely
$ws.Range("A14:H21").EntireRow.Delete()
